$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "On hand" (E) / "Soldered/Recoverable" (F) counts, and the B6 quantity
# tweak, updated to the new component counts. Several of these cells used to
# hold a literal arithmetic formula (e.g. "=3+6"); they are now plain values.

# Row 2 - Ceramic Capacitor, 0.1uF
$ws.Range("E2").Value2 = 0
$ws.Range("F2").Value2 = 8

# Row 3 - Ceramic Capacitor, 10uF
$ws.Range("E3").Value2 = 5

# Row 4 - Ceramic Capacitor, 1uF
$ws.Range("E4").Value2 = 1

# Row 5 - Ceramic Capacitor, 22pF
$ws.Range("E5").Value2 = 13

# Row 6 - Diode, Generic
$ws.Range("B6").Value2 = 22
$ws.Range("E6").Value2 = 144

# Row 7 - Polyfuse
$ws.Range("E7").Value2 = 7

# Row 8 - Resistor, 10k
$ws.Range("E8").Value2 = 12

# Row 9 - Resistor, 22
$ws.Range("E9").Value2 = 12

# Row 10 - Low Profile Tactile Switch
$ws.Range("E10").Value2 = 5

# Row 11 - ATMEGA32U4 Microcontroller
$ws.Range("E11").Value2 = 9

# Row 12 - Molex Mini-B USB C Receptacle
$ws.Range("E12").Value2 = 6

# Row 13 - Crystal, 16MHz
$ws.Range("E13").Value2 = 6

# Row 14 - Stereo Jack 3.5mm
$ws.Range("E14").Value2 = 6

# "Halves to build" modifier used by the H column (=B*$B$18)
$ws.Range("B18").Value2 = 2

# Column I ("Total" = To buy - Qty Needed) is turned into one shared formula
# group spanning I2:I14 (it used to be a separate literal formula per row).
$ws.Range("I2:I14").Formula = "=H2-G2"

# Reflect the author's last selected cell in the saved sheet view.
$ws.Range("A27").Select()
